$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "ActiveSheet: $($ws.Name)"
